# Apply 2022-08-27 daily crime-data update to violent-crime-full-year workbook
# Updates the 2022 (column I) running totals (and a couple of prior-year corrections
# in column E) across the "Citywide Totals", "By Neighborhood" and per-neighborhood
# worksheets.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{Sheet='Citywide Totals'; Cell='I2'; Value=4708},
    @{Sheet='Citywide Totals'; Cell='I3'; Value=4896},
    @{Sheet='Citywide Totals'; Cell='E4'; Value=1967},
    @{Sheet='Citywide Totals'; Cell='I4'; Value=1121},
    @{Sheet='Citywide Totals'; Cell='I5'; Value=451},
    @{Sheet='Citywide Totals'; Cell='I6'; Value=5330},
    @{Sheet='Citywide Totals'; Cell='E7'; Value=25971},
    @{Sheet='Citywide Totals'; Cell='I7'; Value=16506},
    @{Sheet='Uptown'; Cell='I6'; Value=67},
    @{Sheet='Uptown'; Cell='I7'; Value=185},
    @{Sheet='West Ridge'; Cell='I4'; Value=14},
    @{Sheet='West Ridge'; Cell='I7'; Value=174},
    @{Sheet='Bridgeport'; Cell='I4'; Value=3},
    @{Sheet='Bridgeport'; Cell='I5'; Value=1},
    @{Sheet='Bridgeport'; Cell='I6'; Value=31},
    @{Sheet='Bridgeport'; Cell='I7'; Value=92},
    @{Sheet='Grand Crossing'; Cell='I3'; Value=167},
    @{Sheet='Grand Crossing'; Cell='I7'; Value=523},
    @{Sheet='North Lawndale'; Cell='I3'; Value=228},
    @{Sheet='North Lawndale'; Cell='I7'; Value=649},
    @{Sheet='South Deering'; Cell='I2'; Value=51},
    @{Sheet='South Deering'; Cell='I3'; Value=43},
    @{Sheet='South Deering'; Cell='I7'; Value=135},
    @{Sheet='New City'; Cell='I2'; Value=119},
    @{Sheet='New City'; Cell='I3'; Value=108},
    @{Sheet='New City'; Cell='I6'; Value=108},
    @{Sheet='New City'; Cell='I7'; Value=368},
    @{Sheet='By Neighborhood'; Cell='I4'; Value=63},
    @{Sheet='By Neighborhood'; Cell='I5'; Value=52},
    @{Sheet='By Neighborhood'; Cell='I6'; Value=113},
    @{Sheet='By Neighborhood'; Cell='I7'; Value=528},
    @{Sheet='By Neighborhood'; Cell='I8'; Value=1006},
    @{Sheet='By Neighborhood'; Cell='I9'; Value=74},
    @{Sheet='By Neighborhood'; Cell='I11'; Value=251},
    @{Sheet='By Neighborhood'; Cell='I13'; Value=29},
    @{Sheet='By Neighborhood'; Cell='I14'; Value=92},
    @{Sheet='By Neighborhood'; Cell='I19'; Value=463},
    @{Sheet='By Neighborhood'; Cell='I20'; Value=402},
    @{Sheet='By Neighborhood'; Cell='I23'; Value=157},
    @{Sheet='By Neighborhood'; Cell='I29'; Value=1047},
    @{Sheet='By Neighborhood'; Cell='I33'; Value=761},
    @{Sheet='By Neighborhood'; Cell='I36'; Value=223},
    @{Sheet='By Neighborhood'; Cell='I37'; Value=523},
    @{Sheet='By Neighborhood'; Cell='I41'; Value=73},
    @{Sheet='By Neighborhood'; Cell='I42'; Value=565},
    @{Sheet='By Neighborhood'; Cell='I43'; Value=133},
    @{Sheet='By Neighborhood'; Cell='I45'; Value=38},
    @{Sheet='By Neighborhood'; Cell='I46'; Value=33},
    @{Sheet='By Neighborhood'; Cell='I48'; Value=229},
    @{Sheet='By Neighborhood'; Cell='I50'; Value=76},
    @{Sheet='By Neighborhood'; Cell='I51'; Value=174},
    @{Sheet='By Neighborhood'; Cell='I52'; Value=355},
    @{Sheet='By Neighborhood'; Cell='I53'; Value=172},
    @{Sheet='By Neighborhood'; Cell='I54'; Value=367},
    @{Sheet='By Neighborhood'; Cell='I60'; Value=86},
    @{Sheet='By Neighborhood'; Cell='E63'; Value=313},
    @{Sheet='By Neighborhood'; Cell='I63'; Value=59},
    @{Sheet='By Neighborhood'; Cell='I64'; Value=146},
    @{Sheet='By Neighborhood'; Cell='I65'; Value=368},
    @{Sheet='By Neighborhood'; Cell='I67'; Value=649},
    @{Sheet='By Neighborhood'; Cell='I73'; Value=141},
    @{Sheet='By Neighborhood'; Cell='I76'; Value=252},
    @{Sheet='By Neighborhood'; Cell='I78'; Value=232},
    @{Sheet='By Neighborhood'; Cell='I79'; Value=460},
    @{Sheet='By Neighborhood'; Cell='I83'; Value=345},
    @{Sheet='By Neighborhood'; Cell='I84'; Value=135},
    @{Sheet='By Neighborhood'; Cell='I85'; Value=747},
    @{Sheet='By Neighborhood'; Cell='I86'; Value=99},
    @{Sheet='By Neighborhood'; Cell='I88'; Value=149},
    @{Sheet='By Neighborhood'; Cell='I89'; Value=185},
    @{Sheet='By Neighborhood'; Cell='I91'; Value=191},
    @{Sheet='By Neighborhood'; Cell='I93'; Value=98},
    @{Sheet='By Neighborhood'; Cell='I94'; Value=157},
    @{Sheet='By Neighborhood'; Cell='I96'; Value=174},
    @{Sheet='By Neighborhood'; Cell='E101'; Value=25971},
    @{Sheet='By Neighborhood'; Cell='I101'; Value=16506},
    @{Sheet='South Chicago'; Cell='I2'; Value=126},
    @{Sheet='South Chicago'; Cell='I4'; Value=13},
    @{Sheet='South Chicago'; Cell='I5'; Value=14},
    @{Sheet='South Chicago'; Cell='I7'; Value=345},
    @{Sheet='Garfield Park'; Cell='I2'; Value=175},
    @{Sheet='Garfield Park'; Cell='I3'; Value=281},
    @{Sheet='Garfield Park'; Cell='I7'; Value=761},
    @{Sheet='Loop'; Cell='I3'; Value=80},
    @{Sheet='Loop'; Cell='I7'; Value=367},
    @{Sheet='Englewood'; Cell='I2'; Value=307},
    @{Sheet='Englewood'; Cell='I3'; Value=365},
    @{Sheet='Englewood'; Cell='I7'; Value=1047},
    @{Sheet='Chatham'; Cell='I2'; Value=171},
    @{Sheet='Chatham'; Cell='I3'; Value=135},
    @{Sheet='Chatham'; Cell='I7'; Value=463},
    @{Sheet='Lake View'; Cell='I6'; Value=130},
    @{Sheet='Lake View'; Cell='I7'; Value=229},
    @{Sheet='River North'; Cell='I2'; Value=54},
    @{Sheet='River North'; Cell='I7'; Value=252},
    @{Sheet='South Shore'; Cell='I2'; Value=197},
    @{Sheet='South Shore'; Cell='I3'; Value=300},
    @{Sheet='South Shore'; Cell='I6'; Value=183},
    @{Sheet='South Shore'; Cell='I7'; Value=747},
    @{Sheet='Ashburn'; Cell='I3'; Value=32},
    @{Sheet='Ashburn'; Cell='I7'; Value=113},
    @{Sheet='Hermosa'; Cell='I2'; Value=24},
    @{Sheet='Hermosa'; Cell='I7'; Value=73},
    @{Sheet='Humboldt Park'; Cell='I3'; Value=191},
    @{Sheet='Humboldt Park'; Cell='I4'; Value=46},
    @{Sheet='Humboldt Park'; Cell='I6'; Value=160},
    @{Sheet='Humboldt Park'; Cell='I7'; Value=565},
    @{Sheet='Boystown'; Cell='I5'; Value=10},
    @{Sheet='Boystown'; Cell='I6'; Value=29},
    @{Sheet='Rogers Park'; Cell='I2'; Value=54},
    @{Sheet='Rogers Park'; Cell='I3'; Value=60},
    @{Sheet='Rogers Park'; Cell='I7'; Value=232},
    @{Sheet='Jefferson Park'; Cell='I2'; Value=8},
    @{Sheet='Jefferson Park'; Cell='I7'; Value=33},
    @{Sheet='Douglas'; Cell='I2'; Value=44},
    @{Sheet='Douglas'; Cell='I7'; Value=157},
    @{Sheet='Washington Park'; Cell='I3'; Value=71},
    @{Sheet='Washington Park'; Cell='I6'; Value=54},
    @{Sheet='Washington Park'; Cell='I7'; Value=191},
    @{Sheet='Roseland'; Cell='I2'; Value=134},
    @{Sheet='Roseland'; Cell='I6'; Value=132},
    @{Sheet='Roseland'; Cell='I7'; Value=460},
    @{Sheet='Near South Side'; Cell='I2'; Value=41},
    @{Sheet='Near South Side'; Cell='I4'; Value=8},
    @{Sheet='Near South Side'; Cell='I7'; Value=146},
    @{Sheet='Chicago Lawn'; Cell='I2'; Value=111},
    @{Sheet='Chicago Lawn'; Cell='I7'; Value=402},
    @{Sheet='Grand Boulevard'; Cell='I2'; Value=69},
    @{Sheet='Grand Boulevard'; Cell='I7'; Value=223},
    @{Sheet='West Lawn'; Cell='I3'; Value=27},
    @{Sheet='West Lawn'; Cell='I7'; Value=98},
    @{Sheet='Little Village'; Cell='I2'; Value=99},
    @{Sheet='Little Village'; Cell='I7'; Value=355},
    @{Sheet='West Loop'; Cell='I2'; Value=30},
    @{Sheet='West Loop'; Cell='I7'; Value=157},
    @{Sheet='Brighton Park'; Cell='I4'; Value=13},
    @{Sheet='Brighton Park'; Cell='I6'; Value=67},
    @{Sheet='Lincoln Square'; Cell='I2'; Value=20},
    @{Sheet='Lincoln Square'; Cell='I7'; Value=76},
    @{Sheet='Belmont Cragin'; Cell='I4'; Value=21},
    @{Sheet='Belmont Cragin'; Cell='I6'; Value=64},
    @{Sheet='Belmont Cragin'; Cell='I7'; Value=251},
    @{Sheet='Avalon Park'; Cell='I3'; Value=25},
    @{Sheet='Avalon Park'; Cell='I7'; Value=74},
    @{Sheet='Portage Park'; Cell='I3'; Value=47},
    @{Sheet='Portage Park'; Cell='I4'; Value=13},
    @{Sheet='Portage Park'; Cell='I7'; Value=141},
    @{Sheet='United Center'; Cell='I2'; Value=43},
    @{Sheet='United Center'; Cell='I3'; Value=54},
    @{Sheet='United Center'; Cell='I7'; Value=149},
    @{Sheet='Austin'; Cell='I2'; Value=319},
    @{Sheet='Austin'; Cell='I3'; Value=280},
    @{Sheet='Austin'; Cell='I4'; Value=61},
    @{Sheet='Austin'; Cell='I6'; Value=320},
    @{Sheet='Austin'; Cell='I7'; Value=1006},
    @{Sheet='Armour Square'; Cell='I2'; Value=12},
    @{Sheet='Armour Square'; Cell='I7'; Value=52},
    @{Sheet='Streeterville'; Cell='I3'; Value=7},
    @{Sheet='Streeterville'; Cell='I4'; Value=48},
    @{Sheet='Streeterville'; Cell='I7'; Value=99},
    @{Sheet='Little Italy, UIC'; Cell='I2'; Value=35},
    @{Sheet='Little Italy, UIC'; Cell='I6'; Value=70},
    @{Sheet='Little Italy, UIC'; Cell='I7'; Value=174},
    @{Sheet='Morgan Park'; Cell='I2'; Value=29},
    @{Sheet='Morgan Park'; Cell='I7'; Value=86},
    @{Sheet='Hyde Park'; Cell='I2'; Value=27},
    @{Sheet='Hyde Park'; Cell='I6'; Value=75},
    @{Sheet='Hyde Park'; Cell='I7'; Value=133},
    @{Sheet='Logan Square'; Cell='I4'; Value=13},
    @{Sheet='Logan Square'; Cell='I7'; Value=172},
    @{Sheet='Sheffield & DePaul'; Cell='I4'; Value=3},
    @{Sheet='Jackson Park'; Cell='I7'; Value=38},
    @{Sheet='Auburn Gresham'; Cell='I5'; Value=27},
    @{Sheet='Auburn Gresham'; Cell='I6'; Value=135},
    @{Sheet='Auburn Gresham'; Cell='I7'; Value=528},
    @{Sheet='Archer Heights'; Cell='I2'; Value=26},
    @{Sheet='Archer Heights'; Cell='I6'; Value=18},
    @{Sheet='Archer Heights'; Cell='I7'; Value=63}
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}
